$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("K2").Value = 1.791666666666668
$ws.Range("R2").Value = 1.458486584262888
$ws.Range("S2").Value = 1.552746181345467
$ws.Range("K3").Value = 1.791666666666668
$ws.Range("R3").Value = 1.687990434197743
$ws.Range("S3").Value = 1.829652737870454
$ws.Range("K7").Value = 12.93898809523811
$ws.Range("R7").Value = 1.586442583591966
$ws.Range("S7").Value = 1.700608911205746
$ws.Range("K8").Value = 12.93898809523811
$ws.Range("R8").Value = 1.875479296468405
$ws.Range("S8").Value = 2.055366535688525
$ws.Range("K9").Value = 12.93898809523811
$ws.Range("K10").Value = 19.79629629629628
$ws.Range("R10").Value = 1.676945000770297
$ws.Range("S10").Value = 1.806427491177953
$ws.Range("K11").Value = 19.79629629629628
$ws.Range("R11").Value = 2.013021864849877
$ws.Range("S11").Value = 2.224152843076377
$ws.Range("K12").Value = 19.79629629629628
$ws.Range("K13").Value = 19.79629629629628
$ws.Range("R13").Value = 1.676945000770297
$ws.Range("S13").Value = 1.806427491177953
$ws.Range("K14").Value = 19.79629629629628
$ws.Range("R14").Value = 2.013021864849877
$ws.Range("S14").Value = 2.224152843076377
$ws.Range("K15").Value = 12.93898809523811
$ws.Range("K16").Value = 12.93898809523811
$ws.Range("K17").Value = 12.93898809523811
$ws.Range("K18").Value = 12.93898809523811
$ws.Range("K19").Value = 12.93898809523811
$ws.Range("R19").Value = 1.875479296468405
$ws.Range("S19").Value = 2.055366535688525
$ws.Range("K20").Value = 12.93898809523811
$ws.Range("R20").Value = 1.586442583591966
$ws.Range("S20").Value = 1.700608911205746
$ws.Range("K23").Value = 14.47727272727272
$ws.Range("R23").Value = 1.605884483070795
$ws.Range("S23").Value = 1.723253983867794
$ws.Range("K24").Value = 14.47727272727272
$ws.Range("R24").Value = 1.904673198104963
$ws.Range("S24").Value = 2.090962591474336
$ws.Range("K25").Value = 13.76976495726495
$ws.Range("R25").Value = 1.596883662077925
$ws.Range("S25").Value = 1.712764324418727
$ws.Range("K26").Value = 13.76976495726495
$ws.Range("K27").Value = 13.76976495726495
$ws.Range("R27").Value = 1.891133884283326
$ws.Range("S27").Value = 2.074438874323838
$ws.Range("K31").Value = 14.47727272727272
$ws.Range("R31").Value = 1.605884483070795
$ws.Range("S31").Value = 1.723253983867794
$ws.Range("K32").Value = 14.47727272727272
$ws.Range("R32").Value = 1.904673198104963
$ws.Range("S32").Value = 2.090962591474336
$ws.Range("K33").Value = 19.60879629629628
$ws.Range("R33").Value = 1.674333288469303
$ws.Range("S33").Value = 1.803359265239363
$ws.Range("K34").Value = 19.60879629629628
$ws.Range("K35").Value = 19.60879629629628
$ws.Range("R35").Value = 2.008993294560045
$ws.Range("S35").Value = 2.21916989815957
$ws.Range("K38").Value = 0.2777777777777778
$ws.Range("R38").Value = 1.442683896620278
$ws.Range("S38").Value = 1.534625267665953
$ws.Range("K39").Value = 0.2777777777777778
$ws.Range("R39").Value = 1.665380456154149
$ws.Range("S39").Value = 1.802766393442623
$ws.Range("K42").Value = 14.47727272727272
$ws.Range("K43").Value = 14.47727272727272
$ws.Range("R43").Value = 1.605884483070795
$ws.Range("S43").Value = 1.723253983867794
$ws.Range("K44").Value = 14.47727272727272
$ws.Range("R44").Value = 1.904673198104963
$ws.Range("S44").Value = 2.090962591474336
$ws.Range("K50").Value = 12.93898809523811
$ws.Range("K51").Value = 12.93898809523811
$ws.Range("R51").Value = 1.586442583591966
$ws.Range("S51").Value = 1.700608911205746
$ws.Range("K52").Value = 12.93898809523811
$ws.Range("R52").Value = 1.875479296468405
$ws.Range("S52").Value = 2.055366535688525
$ws.Range("K53").Value = 13.76976495726495
$ws.Range("R53").Value = 1.596883662077925
$ws.Range("S53").Value = 1.712764324418727
$ws.Range("K54").Value = 13.76976495726495
$ws.Range("R54").Value = 1.891133884283326
$ws.Range("S54").Value = 2.074438874323838
$ws.Range("K55").Value = 14.47727272727272
$ws.Range("R55").Value = 1.605884483070795
$ws.Range("S55").Value = 1.723253983867794
$ws.Range("K56").Value = 14.47727272727272
$ws.Range("K57").Value = 14.47727272727272
$ws.Range("R57").Value = 1.904673198104963
$ws.Range("S57").Value = 2.090962591474336
$ws.Range("K58").Value = 5.462962962962945
$ws.Range("K59").Value = 5.462962962962945
$ws.Range("R59").Value = 1.745457823876906
$ws.Range("S59").Value = 1.898309736523319
$ws.Range("K60").Value = 5.462962962962945
$ws.Range("R60").Value = 1.49828630419821
$ws.Range("S60").Value = 1.598520446096654
$ws.Range("K61").Value = 19.60879629629628
$ws.Range("R61").Value = 1.674333288469303
$ws.Range("S61").Value = 1.803359265239363
$ws.Range("K62").Value = 19.60879629629628
$ws.Range("R62").Value = 2.008993294560045
$ws.Range("S62").Value = 2.21916989815957
$ws.Range("K63").Value = 19.60879629629628
$ws.Range("K64").Value = 14.47727272727272
$ws.Range("R64").Value = 1.605884483070795
$ws.Range("S64").Value = 1.723253983867794
$ws.Range("K65").Value = 14.47727272727272
$ws.Range("R65").Value = 1.904673198104963
$ws.Range("S65").Value = 2.090962591474336
$ws.Range("K66").Value = 21.28240740740739
$ws.Range("R66").Value = 2.045532758429849
$ws.Range("S66").Value = 2.264453199968389
$ws.Range("K67").Value = 21.28240740740739
$ws.Range("R67").Value = 1.697937058846468
$ws.Range("S67").Value = 1.831120384959332
$ws.Range("K68").Value = 14.47727272727272
$ws.Range("R68").Value = 1.904673198104963
$ws.Range("S68").Value = 2.090962591474336
$ws.Range("K69").Value = 14.47727272727272
$ws.Range("R69").Value = 1.605884483070795
$ws.Range("S69").Value = 1.723253983867794
